# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Most "Price"/"Volume(1h)" cells are plain text in this sheet (the source
# feed formats them like "26.714.80" or "  -1.08%  "), so we just overwrite
# the text in-place. A few Price cells happen to look like plain decimal
# numbers (e.g. "211.74") -- those are entered with a leading apostrophe so
# Excel stores them as text instead of auto-converting to a number, then the
# style is reset back to Normal so no stray text-format styling is left on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.723.21'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '1.599.43'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''211.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = '''19.76'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('D11').Value = '''0.0838'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '1.822.62'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.597.68'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''4.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('E15').Value = '  -2.30%  '
$ws.Range('D16').Value = '''65.12'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '26.698.54'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').Value = '''210.38'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '''6.73'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('E23').Value = '  -1.86%  '
$ws.Range('D24').Value = '''8.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = '''146.86'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = '''7.18'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.01%  '
$ws.Range('D28').Value = '''0.116'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').Value = '''0.0504'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').Value = '''0.672'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.59%  '
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('D35').Value = '1.301.78'
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('D36').Value = '''2.44'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -5.18%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').Value = '''0.844'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').Value = '''5.38'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').Value = '''2.19'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').Value = '''63.98'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').Value = '1.735.51'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '''90.01'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '''0.873'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.60%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '''0.0986'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = '''0.0504'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = '''7.54'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.00%  '
